$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10:E10").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)
$ws.Range("A10:B10").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)
$ws.Range("D10:E10").Copy()
$ws.Range("D12:E12").PasteSpecial(-4122)

$ws.Range("A11:E11").RowHeight = 15.75
$ws.Range("A12:E12").RowHeight = 15.75

$ws.Range("A11").Value = "Розовые вина"
$ws.Range("B11").Value = "Закат Алупки"
$ws.Range("D11").Value = 200
$ws.Range("E11").Value = "konyak_kizilovyi.png"

$ws.Range("A12").Value = "Соки-воды"
$ws.Range("B12").Value = "вода минеральная ""cлезы Кубани"""
$ws.Range("D12").Value = 100
$ws.Range("E12").Value = "chacha.png"

$ws.Range("C11").Value = "Альянико"

$ws.Range("C11").Select()
